# Pricing sheet update: refresh BaneBots wheel link, fix a cost bug (3.5 -> 3.3),
# and add three new line items (M4 Set Screws, M4 Nut, Button caps), pushing the
# PCB / Thank-You-note rows down from 15-16 to 18-19.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# 1) Drop every existing hyperlink so stale Range/Address pairs don't linger
#    once rows move; they get re-created below against the final layout.
$ws.Hyperlinks.Delete()

# 2) Make room: insert three blank rows where row 15 used to be. This pushes
#    the old row15 (PCB) / row16 (Thank You note) down to row18 / row19.
$ws.Rows("15:17").Insert()

# 3) Rewrite every data row (4 through 19) with the post-edit content. Rows
#    2-3 and 1 are untouched by the diff.

# Row 4: BaneBots Wheels - new banebots link, cost formula bug fix 3.5 -> 3.3
$ws.Range("A4").Value2 = "BaneBots Wheels"
$ws.Range("B4").Value2 = "https://banebots.com/t61p-241by/"
$ws.Range("C4").Formula = "=3.3*1.101"
$ws.Range("D4").Value2 = 2
$ws.Range("E4").Formula = "=C4*D4"

# Row 5: Motor (L Shape)
$ws.Range("A5").Value2 = "Motor (L Shape)"
$ws.Range("B5").Value2 = "https://www.aliexpress.us/item/2251832479258088.html"
$ws.Range("C5").Formula = "=15.09*1.101"
$ws.Range("D5").Formula = "=2"
$ws.Range("E5").Formula = "=C5*D5"

# Row 6: Bearing
$ws.Range("A6").Value2 = "Bearing (FMR128ZZ 8x12x3.5mm)"
$ws.Range("B6").Value2 = "https://www.aliexpress.us/item/3256804369338966.html"
$ws.Range("C6").Formula = "=2.29*1.101"
$ws.Range("D6").Value2 = 2
$ws.Range("E6").Formula = "=C6*D6"
$ws.Range("F6").Value2 = "2 out of a pack of 4"

# Row 7: Motor wire
$ws.Range("A7").Value2 = "Motor wire (PH2.0MM, 100MM, 2P, Forward double head)"
$ws.Range("B7").Value2 = "https://www.aliexpress.us/item/3256806095197532.html"
$ws.Range("C7").Formula = "=0.72*1.101"
$ws.Range("D7").Value2 = 2
$ws.Range("E7").Formula = "=C7*D7"

# Row 8: Encoder wire
$ws.Range("A8").Value2 = "Encoder wire (PH2.0MM, 100MM, 4P, Forward double head)"
$ws.Range("B8").Value2 = "https://www.aliexpress.us/item/3256806095197532.html"
$ws.Range("C8").Formula = "=0.93*1.101"
$ws.Range("D8").Value2 = 2
$ws.Range("E8").Formula = "=C8*D8"

# Row 9: Caster wheel
$ws.Range("A9").Value2 = "Caster wheel"
$ws.Range("B9").Value2 = "https://www.adafruit.com/product/3948"
$ws.Range("C9").Formula = "=1.95*1.101"
$ws.Range("D9").Value2 = 2
$ws.Range("E9").Formula = "=C9*D9"

# Row 10: Dowel
$ws.Range("A10").Value2 = "Dowel"
$ws.Range("B10").Value2 = "https://www.amazon.com/dp/B08XQQ69WD"
$ws.Range("C10").Formula = "=4.99*1.101"
$ws.Range("D10").Formula = "=1/25"
$ws.Range("E10").Formula = "=C10*D10"
$ws.Range("F10").Value2 = "1 out of a pack of 25"

# Row 11: M3x6mm screws
$ws.Range("A11").Value2 = "M3x6mm screws"
$ws.Range("B11").Value2 = "https://www.amazon.com/dp/B012TE12CY"
$ws.Range("C11").Formula = "=6.4*1.101"
$ws.Range("D11").Formula = "=9/100"
$ws.Range("E11").Formula = "=C11*D11"
$ws.Range("F11").Value2 = "Need 9/100 screws"

# Row 12: M4x6mm screws
$ws.Range("A12").Value2 = "M4x6mm screws"
$ws.Range("B12").Value2 = "https://www.amazon.com/dp/B07MF33MRJ"
$ws.Range("C12").Formula = "=6.19*1.101"
$ws.Range("D12").Formula = "=7/30"
$ws.Range("E12").Formula = "=C12*D12"
$ws.Range("F12").Value2 = "Need 7/30 screws"

# Row 13: M3 washer
$ws.Range("A13").Value2 = "M3 washer"
$ws.Range("B13").Value2 = "https://www.amazon.com/dp/B0BGH5Y5LQ"
$ws.Range("C13").Formula = "=1.101*5.99"
$ws.Range("D13").Formula = "=3/100"
$ws.Range("E13").Formula = "=C13*D13"
$ws.Range("F13").Value2 = "Need 3/100 washers"

# Row 14: Rechargeable NiMH AA batteries
$ws.Range("A14").Value2 = "Rechargeable NiMH AA batteries"
$ws.Range("B14").Value2 = "https://www.amazon.com/dp/B00CWNMV4G"
$ws.Range("C14").Formula = "=12.34*1.101"
$ws.Range("D14").Value2 = 0
$ws.Range("E14").Formula = "=C14*D14"
$ws.Range("F14").Value2 = "Need 4/8, NOT INCLUDED WITH KIT"

# Row 15 (new): M4 Set Screws
$ws.Range("A15").Value2 = "M4 Set Screws"
$ws.Range("B15").Value2 = "https://www.amazon.com/dp/B01N76NKU6"
$ws.Range("C15").Formula = "=7.09*1.101"
$ws.Range("D15").Formula = "=2/50"
$ws.Range("E15").Formula = "=C15*D15"
$ws.Range("F15").Value2 = "Come with the pulleys I ordered for EV"
$ws.Rows(15).RowHeight = 17

# Row 16 (new): M4 Nut
$ws.Range("A16").Value2 = "M4 Nut"
$ws.Range("B16").Value2 = "https://www.amazon.com/dp/B0BLBLM2BQ/"
$ws.Range("C16").Formula = "=5.99*1.101"
$ws.Range("D16").Formula = "=4/100"
$ws.Range("E16").Formula = "=C16*D16"

# Row 17 (new): Button caps
$ws.Range("A17").Value2 = "Button caps"
$ws.Range("B17").Value2 = "https://www.aliexpress.us/item/2251832666419248.html"
$ws.Range("C17").Formula = "=(0.94*2 + 3.31)*1.101"
$ws.Range("D17").Formula = "=1/50"
$ws.Range("E17").Formula = "=C17*D17"

# Row 18 (was row 15): PCB
$ws.Range("A18").Value2 = "PCB"
$ws.Range("B18").Value2 = "https://jlcpcb.com"
$ws.Range("C18").Formula = "=206.29-6"
$ws.Range("D18").Formula = "=1/10"
$ws.Range("E18").Formula = "=C18*D18"
$ws.Range("F18").Value2 = "Getting 1 of the 10 in the JLCPCB order"

# Row 19 (was row 16): Thank You note
$ws.Range("A19").Value2 = "Thank You note"
$ws.Range("B19").Value2 = "N/A"
$ws.Range("C19").Value2 = 0
$ws.Range("D19").Value2 = 1
$ws.Range("E19").Formula = "=C19*D19"

# 4) Re-create hyperlinks for every Link cell (B2:B18) against the final layout.
$ws.Hyperlinks.Add($ws.Range("B2"), $ws.Range("B2").Value2)
$ws.Hyperlinks.Add($ws.Range("B3"), $ws.Range("B3").Value2)
$ws.Hyperlinks.Add($ws.Range("B4"), $ws.Range("B4").Value2)
$ws.Hyperlinks.Add($ws.Range("B5"), $ws.Range("B5").Value2)
$ws.Hyperlinks.Add($ws.Range("B6"), $ws.Range("B6").Value2)
$ws.Hyperlinks.Add($ws.Range("B7"), $ws.Range("B7").Value2)
$ws.Hyperlinks.Add($ws.Range("B8"), $ws.Range("B8").Value2)
$ws.Hyperlinks.Add($ws.Range("B9"), $ws.Range("B9").Value2)
$ws.Hyperlinks.Add($ws.Range("B10"), $ws.Range("B10").Value2)
$ws.Hyperlinks.Add($ws.Range("B11"), $ws.Range("B11").Value2)
$ws.Hyperlinks.Add($ws.Range("B12"), $ws.Range("B12").Value2)
$ws.Hyperlinks.Add($ws.Range("B13"), $ws.Range("B13").Value2)
$ws.Hyperlinks.Add($ws.Range("B14"), $ws.Range("B14").Value2)
$ws.Hyperlinks.Add($ws.Range("B15"), $ws.Range("B15").Value2)
$ws.Hyperlinks.Add($ws.Range("B16"), $ws.Range("B16").Value2)
$ws.Hyperlinks.Add($ws.Range("B17"), $ws.Range("B17").Value2)
$ws.Hyperlinks.Add($ws.Range("B18"), $ws.Range("B18").Value2)

# 5) Selection matches the author's last cursor position.
$ws.Range("F17").Select()

$wb.Application.Calculate()
